$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 22 (2025-09 stats) per commit "Update stats for 2025-09"
$ws.Range("B22").Value = 6293
$ws.Range("D22").Value = 5845608
$ws.Range("E22").Value = 928.9064039408867
$ws.Range("F22").Value = 8.331898777758639
$ws.Range("H22").Value = 27.12375716559816
